$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14's phone number was stored as text ("71277628"); the redemption
# recorded here should use a real numeric value, matching the other rows.
$ws.Range("A14").Value = 71277628

# Append the new redemption: phone 71277620 redeems 76.0 points.
# Format the phone column as text first so the long numeric string isn't
# auto-converted to a number (consistent with how the other phone-as-text
# values in this sheet are stored).
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "71277620"
$ws.Range("B15").Value = 76
$ws.Range("C15").Value = "2025-08-18T17:04:15"
